# PDP11Sim_TestPlan.xlsx - mark the newly-added Conditional Branch test
# cases (CB_01 .. CB_15, rows 16-30) as implemented in the "Implemented?"
# column, matching the "X" marker already used on rows 8-14.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B16:B30").Value = "X"

# Leave the view where the author was last working in the sheet.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 16
$ws.Range("A29").Select()

$wb.Save()
